# Auto-generated Excel COM-interop script
# Updates market-price columns (H-N) on several leve rows across all 8 sheets
# per the scheduled market-data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

$ws.Range("H70").Value = 4582.476
$ws.Range("J70").Value = 4582.476
$ws.Range("L70").Value = 13747.428
$ws.Range("N70").Value = -14287.428

$ws.Range("H73").Value = 4582.476
$ws.Range("J73").Value = 4582.476
$ws.Range("L73").Value = 13747.428
$ws.Range("N73").Value = -15619.428

$ws.Range("H100").Value = 597.5
$ws.Range("I100").Value = 597.5
$ws.Range("K100").Value = 597.5
$ws.Range("M100").Value = -56.5

$ws.Range("H107").Value = 570.25
$ws.Range("I107").Value = 627.5714
$ws.Range("K107").Value = 627.5714
$ws.Range("M107").Value = 1292.4286

$ws.Range("H111").Value = 1699.8
$ws.Range("I111").Value = 833.3333
$ws.Range("J111").Value = 2999.5
$ws.Range("K111").Value = 2499.9999
$ws.Range("L111").Value = 8998.5
$ws.Range("M111").Value = 567.0001000000002
$ws.Range("N111").Value = -15132.5

$ws.Range("H113").Value = 4765.857
$ws.Range("I113").Value = 4499.5
$ws.Range("J113").Value = 4872.4
$ws.Range("K113").Value = 4499.5
$ws.Range("L113").Value = 4872.4
$ws.Range("M113").Value = -1245.5
$ws.Range("N113").Value = -11380.4

$ws.Range("H118").Value = 1221.091
$ws.Range("I118").Value = 937.2222
$ws.Range("J118").Value = 2498.5
$ws.Range("K118").Value = 2811.6666
$ws.Range("L118").Value = 7495.5
$ws.Range("M118").Value = -1154.6666
$ws.Range("N118").Value = -10809.5

$ws.Range("H135").Value = 1133.375
$ws.Range("I135").Value = 869.6667
$ws.Range("K135").Value = 7827.0003
$ws.Range("M135").Value = -5292.0003

$ws.Range("H141").Value = 3615.2727
$ws.Range("I141").Value = 3226.6667
$ws.Range("K141").Value = 9680.000100000001
$ws.Range("M141").Value = -4500.000100000001


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H45").Value = 1947.25
$ws.Range("I45").Value = 1947.25
$ws.Range("K45").Value = 1947.25
$ws.Range("M45").Value = -1570.25

$ws.Range("H55").Value = 26749.625
$ws.Range("J55").Value = 28333.334
$ws.Range("L55").Value = 28333.334
$ws.Range("N55").Value = -28963.334

$ws.Range("H61").Value = 1936.25
$ws.Range("I61").Value = 1936.25
$ws.Range("K61").Value = 1936.25
$ws.Range("M61").Value = -1724.25

$ws.Range("H94").Value = 86775.664
$ws.Range("J94").Value = 86775.664
$ws.Range("L94").Value = 86775.664
$ws.Range("N94").Value = -88577.664

$ws.Range("H132").Value = 2454.3572
$ws.Range("I132").Value = 2454.3572
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7363.071599999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4833.071599999999
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1936.25
$ws.Range("I136").Value = 1936.25
$ws.Range("K136").Value = 5808.75
$ws.Range("M136").Value = -3258.75


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2367.611
$ws.Range("I86").Value = 1375.9333
$ws.Range("J86").Value = 7326
$ws.Range("K86").Value = 1375.9333
$ws.Range("L86").Value = 7326
$ws.Range("M86").Value = -252.9332999999999
$ws.Range("N86").Value = -9572

$ws.Range("H89").Value = 2367.611
$ws.Range("I89").Value = 1375.9333
$ws.Range("J89").Value = 7326
$ws.Range("K89").Value = 6879.666499999999
$ws.Range("L89").Value = 36630
$ws.Range("M89").Value = -1263.666499999999
$ws.Range("N89").Value = -47862

$ws.Range("H99").Value = 2318
$ws.Range("I99").Value = 2066.6667
$ws.Range("J99").Value = 2695
$ws.Range("K99").Value = 2066.6667
$ws.Range("L99").Value = 2695
$ws.Range("M99").Value = -568.6667000000002
$ws.Range("N99").Value = -5691


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 85
$ws.Range("I3").Value = 85
$ws.Range("K3").Value = 85
$ws.Range("M3").Value = 28

$ws.Range("H31").Value = 3226.8
$ws.Range("I31").Value = 2810
$ws.Range("J31").Value = 3504.6667
$ws.Range("K31").Value = 2810
$ws.Range("L31").Value = 3504.6667
$ws.Range("M31").Value = -2515
$ws.Range("N31").Value = -4094.6667

$ws.Range("H34").Value = 3226.8
$ws.Range("I34").Value = 2810
$ws.Range("J34").Value = 3504.6667
$ws.Range("K34").Value = 2810
$ws.Range("L34").Value = 3504.6667
$ws.Range("M34").Value = -2608
$ws.Range("N34").Value = -3908.6667

$ws.Range("H132").Value = 1695.8334
$ws.Range("J132").Value = 1999.5
$ws.Range("L132").Value = 5998.5
$ws.Range("N132").Value = -11058.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 312.5
$ws.Range("I2").Value = 34.18182
$ws.Range("K2").Value = 205.09092
$ws.Range("M2").Value = -92.09092000000001

$ws.Range("H9").Value = 99
$ws.Range("I9").Value = 99
$ws.Range("K9").Value = 297
$ws.Range("M9").Value = -73

$ws.Range("H38").Value = 35.2
$ws.Range("I38").Value = 32.5
$ws.Range("J38").Value = 37
$ws.Range("K38").Value = 97.5
$ws.Range("L38").Value = 111
$ws.Range("M38").Value = 249.5
$ws.Range("N38").Value = -805

$ws.Range("H46").Value = 1772
$ws.Range("I46").Value = 696
$ws.Range("K46").Value = 2088
$ws.Range("M46").Value = -1997

$ws.Range("H97").Value = 898.9091
$ws.Range("I97").Value = 1384.6
$ws.Range("J97").Value = 494.16666
$ws.Range("K97").Value = 4153.799999999999
$ws.Range("L97").Value = 1482.49998
$ws.Range("M97").Value = -3657.799999999999
$ws.Range("N97").Value = -2474.49998

$ws.Range("H113").Value = 579.1579
$ws.Range("J113").Value = 605.8461
$ws.Range("L113").Value = 1817.5383
$ws.Range("N113").Value = -6157.5383

$ws.Range("H116").Value = 3448
$ws.Range("I116").Value = 2900
$ws.Range("K116").Value = 8700
$ws.Range("M116").Value = -5258

$ws.Range("H117").Value = 650
$ws.Range("I117").Value = 650
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1950
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 1492
$ws.Range("N117").ClearContents()


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 23666.334
$ws.Range("J44").Value = 23666.334
$ws.Range("L44").Value = 23666.334
$ws.Range("N44").Value = -24858.334

$ws.Range("H47").Value = 22500
$ws.Range("J47").Value = 22500
$ws.Range("L47").Value = 22500
$ws.Range("N47").Value = -23636

$ws.Range("H102").Value = 1733.3334
$ws.Range("I102").Value = 1133.3334
$ws.Range("K102").Value = 1133.3334
$ws.Range("M102").Value = 488.6666


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2803.818
$ws.Range("I46").Value = 2540.5881
$ws.Range("J46").Value = 3698.8
$ws.Range("K46").Value = 2540.5881
$ws.Range("L46").Value = 3698.8
$ws.Range("M46").Value = -2352.5881
$ws.Range("N46").Value = -4074.8

$ws.Range("H55").Value = 3024.6667
$ws.Range("I55").Value = 2537.5
$ws.Range("J55").Value = 3999
$ws.Range("K55").Value = 2537.5
$ws.Range("L55").Value = 3999
$ws.Range("M55").Value = -2364.5
$ws.Range("N55").Value = -4345

$ws.Range("H105").Value = 75000
$ws.Range("J105").Value = 75000
$ws.Range("L105").Value = 75000
$ws.Range("N105").Value = -81988

$ws.Range("H122").Value = 3264.1
$ws.Range("I122").Value = 3264.1
$ws.Range("K122").Value = 9792.299999999999
$ws.Range("M122").Value = -7342.299999999999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2527.7778
$ws.Range("I2").Value = 2470.5881
$ws.Range("J2").Value = 3500
$ws.Range("K2").Value = 2470.5881
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = -2358.5881
$ws.Range("N2").Value = -3724

$ws.Range("H4").Value = 35166.668

$ws.Range("H6").Value = 150
$ws.Range("J6").Value = 150
$ws.Range("L6").Value = 150
$ws.Range("N6").Value = -380

$ws.Range("H41").Value = 15813.5
$ws.Range("I41").Value = 16469
$ws.Range("K41").Value = 16469
$ws.Range("M41").Value = -16079

$ws.Range("H45").Value = 8100
$ws.Range("I45").Value = 3800
$ws.Range("J45").Value = 10250
$ws.Range("K45").Value = 3800
$ws.Range("L45").Value = 10250
$ws.Range("M45").Value = -3309
$ws.Range("N45").Value = -11232

$ws.Range("H132").Value = 1482.8334
$ws.Range("I132").Value = 1474.25
$ws.Range("K132").Value = 4422.75
$ws.Range("M132").Value = -1892.75

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 1858.2632
$ws.Range("J136").Value = 1398.5
$ws.Range("L136").Value = 4195.5
$ws.Range("N136").Value = -9295.5

